# Zeitplanung.xlsx - Anpassung des Zeitplans (vh)
#
# Applies the target edit to the "Aufgabenliste Projekt 1" sheet:
#   - marks a few tasks as completed (% erledigt -> 1)
#   - records an actual-completion date for one of them
#   - swaps the person assigned to two tasks and halves their logged hours
#   - leaves the cursor/selection where the editor ended up working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: "Storyboard-Zeichungen" follow-up task -> done -------------
$ws.Range("E20").Value = 1
# F20 is the calculated "Fortschritt" column; it keeps its table formula
# and simply recalculates to the new percentage.

# --- Row 21: task finished, with an actual completion date --------------
$ws.Range("E21").Value = 1
# This row's Fortschritt cell loses its live formula and becomes the
# literal value that was on screen when the row was marked done.
$ws.Range("F21").Value = 1
# Pick up the date format already used by the other "tatsächliche
# Fertigstellung" entries in this column, then stamp the date itself.
$ws.Range("H15").Copy()
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("H21").Value = 42341

# --- Row 22: another task -> done ---------------------------------------
$ws.Range("E22").Value = 1

# --- Rows 36/37: re-assign who does which scene and correct hours -------
# "(S) - Schrank, Mobile, Fenster-Szene fertig stellen" goes to Tobias,
# "(S) - Bett + Kommode + Junge Szene fertig stellen" goes to Viktoria;
# both tasks' logged effort is corrected from 4h down to 2h.
$ws.Range("D36").Value = "Tobias"
$ws.Range("G36").Value = "2h"
$ws.Range("D37").Value = "Viktoria"
$ws.Range("G37").Value = "2h "

# --- Leave the view where the editor ended up -----------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 37
$win.ScrollColumn = 2
$ws.Range("D38").Select()
